$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell reference -> new text value, taken from the
# refreshed cryptos feed (price + 1h volume % columns), plus the two rows
# whose coin/link swapped position (HuobiToken <-> ARBITRUM).
$updates = @(
    @{ Cell = 'D2'; Value = '28.922.91' }
    @{ Cell = 'E2'; Value = '  -2.58%  ' }
    @{ Cell = 'D3'; Value = '1.879.66' }
    @{ Cell = 'E3'; Value = '  -3.50%  ' }
    @{ Cell = 'D4'; Value = '1.002' }
    @{ Cell = 'E4'; Value = '  +0.66%  ' }
    @{ Cell = 'D5'; Value = '330.91' }
    @{ Cell = 'E5'; Value = '  -3.36%  ' }
    @{ Cell = 'D6'; Value = '1.002' }
    @{ Cell = 'E6'; Value = '  +0.67%  ' }
    @{ Cell = 'D7'; Value = '0.4599' }
    @{ Cell = 'E7'; Value = '  -4.16%  ' }
    @{ Cell = 'D8'; Value = '0.4064' }
    @{ Cell = 'E8'; Value = '  -1.73%  ' }
    @{ Cell = 'D9'; Value = '47.87' }
    @{ Cell = 'E9'; Value = '  -2.04%  ' }
    @{ Cell = 'D10'; Value = '0.07970' }
    @{ Cell = 'E10'; Value = '  -3.73%  ' }
    @{ Cell = 'D11'; Value = '0.9883' }
    @{ Cell = 'E11'; Value = '  -5.89%  ' }
    @{ Cell = 'D12'; Value = '21.59' }
    @{ Cell = 'E12'; Value = '  -4.86%  ' }
    @{ Cell = 'D13'; Value = '1.912.04' }
    @{ Cell = 'E13'; Value = '  -0.64%  ' }
    @{ Cell = 'D14'; Value = '5.893' }
    @{ Cell = 'E14'; Value = '  -4.30%  ' }
    @{ Cell = 'D15'; Value = '7.058' }
    @{ Cell = 'E15'; Value = '  -5.33%  ' }
    @{ Cell = 'D16'; Value = '1.004' }
    @{ Cell = 'E16'; Value = '  +0.84%  ' }
    @{ Cell = 'D17'; Value = '88.24' }
    @{ Cell = 'E17'; Value = '  -5.26%  ' }
    @{ Cell = 'D18'; Value = '0.00001029' }
    @{ Cell = 'E18'; Value = '  -3.83%  ' }
    @{ Cell = 'D19'; Value = '0.06579' }
    @{ Cell = 'E19'; Value = '  -1.38%  ' }
    @{ Cell = 'D20'; Value = '17.38' }
    @{ Cell = 'E20'; Value = '  -4.05%  ' }
    @{ Cell = 'D21'; Value = '1.004' }
    @{ Cell = 'E21'; Value = '  +0.87%  ' }
    @{ Cell = 'D22'; Value = '28.879.82' }
    @{ Cell = 'E22'; Value = '  -2.60%  ' }
    @{ Cell = 'D23'; Value = '5.403' }
    @{ Cell = 'E23'; Value = '  -3.80%  ' }
    @{ Cell = 'D24'; Value = '11.48' }
    @{ Cell = 'E24'; Value = '  +1.56%  ' }
    @{ Cell = 'D25'; Value = '2.205' }
    @{ Cell = 'E25'; Value = '  -2.35%  ' }
    @{ Cell = 'D26'; Value = '2.099.11' }
    @{ Cell = 'E26'; Value = '  -2.54%  ' }
    @{ Cell = 'D27'; Value = '156.73' }
    @{ Cell = 'E27'; Value = '  -2.80%  ' }
    @{ Cell = 'D28'; Value = '19.50' }
    @{ Cell = 'E28'; Value = '  -3.71%  ' }
    @{ Cell = 'D29'; Value = '2.078' }
    @{ Cell = 'E29'; Value = '  -6.07%  ' }
    @{ Cell = 'D30'; Value = '5.466' }
    @{ Cell = 'E30'; Value = '  -3.00%  ' }
    @{ Cell = 'D31'; Value = '117.29' }
    @{ Cell = 'E31'; Value = '  -4.26%  ' }
    @{ Cell = 'D32'; Value = '1.015' }
    @{ Cell = 'E32'; Value = '  -0.97%  ' }
    @{ Cell = 'D33'; Value = '0.09327' }
    @{ Cell = 'E33'; Value = '  -3.48%  ' }
    @{ Cell = 'B34'; Value = 'ARBITRUM' }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' }
    @{ Cell = 'D34'; Value = '1.397' }
    @{ Cell = 'E34'; Value = '  -5.13%  ' }
    @{ Cell = 'B35'; Value = 'HuobiToken' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Cell = 'D35'; Value = '3.473' }
    @{ Cell = 'E35'; Value = '  -5.14%  ' }
    @{ Cell = 'D36'; Value = '5.268' }
    @{ Cell = 'E36'; Value = '  -4.33%  ' }
    @{ Cell = 'D37'; Value = '0.06033' }
    @{ Cell = 'E37'; Value = '  -3.27%  ' }
    @{ Cell = 'D38'; Value = '0.02225' }
    @{ Cell = 'E38'; Value = '  -3.71%  ' }
    @{ Cell = 'D39'; Value = '8.287' }
    @{ Cell = 'E39'; Value = '  -4.88%  ' }
    @{ Cell = 'D40'; Value = '1.169' }
    @{ Cell = 'E40'; Value = '  -3.10%  ' }
    @{ Cell = 'D41'; Value = '1.002' }
    @{ Cell = 'E41'; Value = '  +0.87%  ' }
    @{ Cell = 'D42'; Value = '0.5763' }
    @{ Cell = 'E42'; Value = '  -5.89%  ' }
    @{ Cell = 'D43'; Value = '0.1825' }
    @{ Cell = 'E43'; Value = '  -4.65%  ' }
    @{ Cell = 'D44'; Value = '10.07' }
    @{ Cell = 'E44'; Value = '  -5.84%  ' }
    @{ Cell = 'D45'; Value = '1.241' }
    @{ Cell = 'E45'; Value = '  -2.30%  ' }
    @{ Cell = 'D46'; Value = '0.07467' }
    @{ Cell = 'E46'; Value = '  +2.97%  ' }
    @{ Cell = 'D47'; Value = '2.271' }
    @{ Cell = 'E47'; Value = '  -1.87%  ' }
    @{ Cell = 'D48'; Value = '11.97' }
    @{ Cell = 'E48'; Value = '  -5.38%  ' }
    @{ Cell = 'D49'; Value = '0.5439' }
    @{ Cell = 'E49'; Value = '  -5.12%  ' }
    @{ Cell = 'D50'; Value = '1.899' }
    @{ Cell = 'E50'; Value = '  -5.39%  ' }
    @{ Cell = 'D51'; Value = '111.00' }
    @{ Cell = 'E51'; Value = '  -2.36%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Cell.Substring(0,1) -eq "D") {
        # Column D holds price text that frequently looks like a plain
        # decimal number ("1.002", "330.91", ...). Excel auto-converts such
        # strings to a Number when assigned directly, which would change the
        # cell type away from the original text cell. Forcing the cell to
        # the Text number format before the write keeps it a string, and
        # ClearFormats() afterwards drops the temporary "@" format again so
        # the cell style index is left exactly as it was.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.ClearFormats()
    } else {
        $rng.Value = $u.Value
    }
}
